$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D2 currently holds "ipfdb.database.windows.net,1433" (the non-default port
# variant). Change it to match the rest of the Host column ("ipfdb.database.windows.net").
$ws.Range("D2").Value = "ipfdb.database.windows.net"

# Reflect the new active selection left in the saved worksheet (the cell the
# user edited).
$ws.Range("D2").Select()
